$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the existing "index" column cell (A52) down into the new rows (A53:A57)
$ws.Range("A52").Copy()
$ws.Range("A53:A57").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# B (lamda_1) and C (lamda_2) are constant across every data row (2..57) -- update them in one shot
$ws.Range("B2:B57").Value2 = 33.94444444444444
$ws.Range("C2:C57").Value2 = 1.95

# A (index), D (dic_nbre_clients_poisson_2_keys) and E (dic_nbre_clients_prob_poisson_2_values) per row
$ws.Cells.Item(2, 1).Value2 = 0
$ws.Cells.Item(2, 4).Value2 = 0
$ws.Cells.Item(2, 5).Value2 = 0.142
$ws.Cells.Item(3, 1).Value2 = 1
$ws.Cells.Item(3, 4).Value2 = 3
$ws.Cells.Item(3, 5).Value2 = 0.004
$ws.Cells.Item(4, 1).Value2 = 2
$ws.Cells.Item(4, 4).Value2 = 4
$ws.Cells.Item(4, 5).Value2 = 0.015
$ws.Cells.Item(5, 1).Value2 = 3
$ws.Cells.Item(5, 4).Value2 = 5
$ws.Cells.Item(5, 5).Value2 = 0.019
$ws.Cells.Item(6, 1).Value2 = 4
$ws.Cells.Item(6, 4).Value2 = 6
$ws.Cells.Item(6, 5).Value2 = 0.039
$ws.Cells.Item(7, 1).Value2 = 5
$ws.Cells.Item(7, 4).Value2 = 7
$ws.Cells.Item(7, 5).Value2 = 0.043
$ws.Cells.Item(8, 1).Value2 = 6
$ws.Cells.Item(8, 4).Value2 = 8
$ws.Cells.Item(8, 5).Value2 = 0.044
$ws.Cells.Item(9, 1).Value2 = 7
$ws.Cells.Item(9, 4).Value2 = 9
$ws.Cells.Item(9, 5).Value2 = 0.052
$ws.Cells.Item(10, 1).Value2 = 8
$ws.Cells.Item(10, 4).Value2 = 10
$ws.Cells.Item(10, 5).Value2 = 0.042
$ws.Cells.Item(11, 1).Value2 = 9
$ws.Cells.Item(11, 4).Value2 = 11
$ws.Cells.Item(11, 5).Value2 = 0.025
$ws.Cells.Item(12, 1).Value2 = 10
$ws.Cells.Item(12, 4).Value2 = 12
$ws.Cells.Item(12, 5).Value2 = 0.033
$ws.Cells.Item(13, 1).Value2 = 11
$ws.Cells.Item(13, 4).Value2 = 13
$ws.Cells.Item(13, 5).Value2 = 0.018
$ws.Cells.Item(14, 1).Value2 = 12
$ws.Cells.Item(14, 4).Value2 = 14
$ws.Cells.Item(14, 5).Value2 = 0.03
$ws.Cells.Item(15, 1).Value2 = 13
$ws.Cells.Item(15, 4).Value2 = 15
$ws.Cells.Item(15, 5).Value2 = 0.033
$ws.Cells.Item(16, 1).Value2 = 14
$ws.Cells.Item(16, 4).Value2 = 16
$ws.Cells.Item(16, 5).Value2 = 0.031
$ws.Cells.Item(17, 1).Value2 = 15
$ws.Cells.Item(17, 4).Value2 = 17
$ws.Cells.Item(17, 5).Value2 = 0.046
$ws.Cells.Item(18, 1).Value2 = 16
$ws.Cells.Item(18, 4).Value2 = 18
$ws.Cells.Item(18, 5).Value2 = 0.027
$ws.Cells.Item(19, 1).Value2 = 17
$ws.Cells.Item(19, 4).Value2 = 19
$ws.Cells.Item(19, 5).Value2 = 0.029
$ws.Cells.Item(20, 1).Value2 = 18
$ws.Cells.Item(20, 4).Value2 = 20
$ws.Cells.Item(20, 5).Value2 = 0.022
$ws.Cells.Item(21, 1).Value2 = 19
$ws.Cells.Item(21, 4).Value2 = 21
$ws.Cells.Item(21, 5).Value2 = 0.026
$ws.Cells.Item(22, 1).Value2 = 20
$ws.Cells.Item(22, 4).Value2 = 22
$ws.Cells.Item(22, 5).Value2 = 0.021
$ws.Cells.Item(23, 1).Value2 = 21
$ws.Cells.Item(23, 4).Value2 = 23
$ws.Cells.Item(23, 5).Value2 = 0.034
$ws.Cells.Item(24, 1).Value2 = 22
$ws.Cells.Item(24, 4).Value2 = 24
$ws.Cells.Item(24, 5).Value2 = 0.02
$ws.Cells.Item(25, 1).Value2 = 23
$ws.Cells.Item(25, 4).Value2 = 25
$ws.Cells.Item(25, 5).Value2 = 0.024
$ws.Cells.Item(26, 1).Value2 = 24
$ws.Cells.Item(26, 4).Value2 = 26
$ws.Cells.Item(26, 5).Value2 = 0.016
$ws.Cells.Item(27, 1).Value2 = 25
$ws.Cells.Item(27, 4).Value2 = 27
$ws.Cells.Item(27, 5).Value2 = 0.02
$ws.Cells.Item(28, 1).Value2 = 26
$ws.Cells.Item(28, 4).Value2 = 28
$ws.Cells.Item(28, 5).Value2 = 0.013
$ws.Cells.Item(29, 1).Value2 = 27
$ws.Cells.Item(29, 4).Value2 = 29
$ws.Cells.Item(29, 5).Value2 = 0.018
$ws.Cells.Item(30, 1).Value2 = 28
$ws.Cells.Item(30, 4).Value2 = 30
$ws.Cells.Item(30, 5).Value2 = 0.006
$ws.Cells.Item(31, 1).Value2 = 29
$ws.Cells.Item(31, 4).Value2 = 31
$ws.Cells.Item(31, 5).Value2 = 0.011
$ws.Cells.Item(32, 1).Value2 = 30
$ws.Cells.Item(32, 4).Value2 = 32
$ws.Cells.Item(32, 5).Value2 = 0.008
$ws.Cells.Item(33, 1).Value2 = 31
$ws.Cells.Item(33, 4).Value2 = 33
$ws.Cells.Item(33, 5).Value2 = 0.011
$ws.Cells.Item(34, 1).Value2 = 32
$ws.Cells.Item(34, 4).Value2 = 34
$ws.Cells.Item(34, 5).Value2 = 0.01
$ws.Cells.Item(35, 1).Value2 = 33
$ws.Cells.Item(35, 4).Value2 = 35
$ws.Cells.Item(35, 5).Value2 = 0.009000000000000001
$ws.Cells.Item(36, 1).Value2 = 34
$ws.Cells.Item(36, 4).Value2 = 36
$ws.Cells.Item(36, 5).Value2 = 0.008
$ws.Cells.Item(37, 1).Value2 = 35
$ws.Cells.Item(37, 4).Value2 = 37
$ws.Cells.Item(37, 5).Value2 = 0.001
$ws.Cells.Item(38, 1).Value2 = 36
$ws.Cells.Item(38, 4).Value2 = 38
$ws.Cells.Item(38, 5).Value2 = 0.005
$ws.Cells.Item(39, 1).Value2 = 37
$ws.Cells.Item(39, 4).Value2 = 39
$ws.Cells.Item(39, 5).Value2 = 0.009000000000000001
$ws.Cells.Item(40, 1).Value2 = 38
$ws.Cells.Item(40, 4).Value2 = 40
$ws.Cells.Item(40, 5).Value2 = 0.006
$ws.Cells.Item(41, 1).Value2 = 39
$ws.Cells.Item(41, 4).Value2 = 41
$ws.Cells.Item(41, 5).Value2 = 0.003
$ws.Cells.Item(42, 1).Value2 = 40
$ws.Cells.Item(42, 4).Value2 = 42
$ws.Cells.Item(42, 5).Value2 = 0.001
$ws.Cells.Item(43, 1).Value2 = 41
$ws.Cells.Item(43, 4).Value2 = 43
$ws.Cells.Item(43, 5).Value2 = 0.002
$ws.Cells.Item(44, 1).Value2 = 42
$ws.Cells.Item(44, 4).Value2 = 44
$ws.Cells.Item(44, 5).Value2 = 0.002
$ws.Cells.Item(45, 1).Value2 = 43
$ws.Cells.Item(45, 4).Value2 = 45
$ws.Cells.Item(45, 5).Value2 = 0.002
$ws.Cells.Item(46, 1).Value2 = 44
$ws.Cells.Item(46, 4).Value2 = 46
$ws.Cells.Item(46, 5).Value2 = 0.001
$ws.Cells.Item(47, 1).Value2 = 45
$ws.Cells.Item(47, 4).Value2 = 47
$ws.Cells.Item(47, 5).Value2 = 0.001
$ws.Cells.Item(48, 1).Value2 = 46
$ws.Cells.Item(48, 4).Value2 = 48
$ws.Cells.Item(48, 5).Value2 = 0.005
$ws.Cells.Item(49, 1).Value2 = 47
$ws.Cells.Item(49, 4).Value2 = 49
$ws.Cells.Item(49, 5).Value2 = 0.003
$ws.Cells.Item(50, 1).Value2 = 48
$ws.Cells.Item(50, 4).Value2 = 51
$ws.Cells.Item(50, 5).Value2 = 0.002
$ws.Cells.Item(51, 1).Value2 = 49
$ws.Cells.Item(51, 4).Value2 = 52
$ws.Cells.Item(51, 5).Value2 = 0.001
$ws.Cells.Item(52, 1).Value2 = 50
$ws.Cells.Item(52, 4).Value2 = 56
$ws.Cells.Item(52, 5).Value2 = 0.001
$ws.Cells.Item(53, 1).Value2 = 51
$ws.Cells.Item(53, 4).Value2 = 58
$ws.Cells.Item(53, 5).Value2 = 0.001
$ws.Cells.Item(54, 1).Value2 = 52
$ws.Cells.Item(54, 4).Value2 = 59
$ws.Cells.Item(54, 5).Value2 = 0.001
$ws.Cells.Item(55, 1).Value2 = 53
$ws.Cells.Item(55, 4).Value2 = 62
$ws.Cells.Item(55, 5).Value2 = 0.001
$ws.Cells.Item(56, 1).Value2 = 54
$ws.Cells.Item(56, 4).Value2 = 64
$ws.Cells.Item(56, 5).Value2 = 0.001
$ws.Cells.Item(57, 1).Value2 = 55
$ws.Cells.Item(57, 4).Value2 = 68
$ws.Cells.Item(57, 5).Value2 = 0.001
